# Auto-generated edit script for LOQ4037.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet shrinks from 25 to 24 rows - remove the trailing row 25
$ws.Rows.Item(25).Delete()

# Clear cells that should become empty
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# Set final cell values for rows 1-24
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOQ4037"
$ws.Range("C2").Value = "LOQ4037"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Química Orgânica I"
$ws.Range("C3").Value = " Química Orgânica I"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Organic Chemistry I"
$ws.Range("C4").Value = "Organic Chemistry I"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2020"
$ws.Range("C8").Value = "01/01/2020"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EQD-3,EQN-4"
$ws.Range("C9").Value = "EQD-3,EQN-4"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "210064 - Eduardo Rezende Triboni"
$ws.Range("C10").Value = "210064 - Eduardo Rezende Triboni"

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Overview - Introduce and teach concepts of organic chemistry as important tools for understanding strategies and industrial and technological operations. Address social and environmental issues with which chemical engineering is related, making them thus able to exercise Chemical Engineer function, and realize the changes that are necessary.Specific - Understand and describe the mechanism of organic reactions and their importance to the improvement and development of synthetic manufacturing processes and formulation stages. Deepening the concept of structure-reactivity and properties of materials."
$ws.Range("C11").Value = "Overview - Introduce and teach concepts of organic chemistry as important tools for understanding strategies and industrial and technological operations. Address social and environmental issues with which chemical engineering is related, making them thus able to exercise Chemical Engineer function, and realize the changes that are necessary.Specific - Understand and describe the mechanism of organic reactions and their importance to the improvement and development of synthetic manufacturing processes and formulation stages. Deepening the concept of structure-reactivity and properties of materials."

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "General property of organic compounds. Physical properties, reactions of aliphatic and aromatic hydrocarbons, organic halides, ethers, alcohols and structural characteristics as stereochemistry and structure-reactivity."
$ws.Range("C14").Value = "General property of organic compounds. Physical properties, reactions of aliphatic and aromatic hydrocarbons, organic halides, ethers, alcohols and structural characteristics as stereochemistry and structure-reactivity."

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1.Bronsted and Lewis acid of the organic compounds2.Alkanes - obtaining processes, physical properties, conformational analysis. Radical substitution reaction.3.Constitutional isomerism and Stereochemistry. Chirality, nomenclature R/S. Polarimeter and Techniques for resolution of stereoisomers.4.Alkyl halides - Nucleophilic Substitution, SN1, SN2, E1, E2.5.Alkenes, alkadienes and alkynes - Physical and chemical properties. Electrophilic addition reaction (hidrohalogenation, hydration, halogenation, Diels-Alder, reduction and oxidation). Conjugated Addition in dienes (thermodynamic and kinetic product).6 Background of NMR, InfraRed, UV and Fluorescence techniques7.Aromatic compounds - Physical properties. Aromatic Eletrophilic Substitution . Effect of Substituent Groups. Aromatic Nucleophilic Substitution.8.Alcohols and ethers - physical properties, reactions and mechanisms."
$ws.Range("C16").Value = "1.Bronsted and Lewis acid of the organic compounds2.Alkanes - obtaining processes, physical properties, conformational analysis. Radical substitution reaction.3.Constitutional isomerism and Stereochemistry. Chirality, nomenclature R/S. Polarimeter and Techniques for resolution of stereoisomers.4.Alkyl halides - Nucleophilic Substitution, SN1, SN2, E1, E2.5.Alkenes, alkadienes and alkynes - Physical and chemical properties. Electrophilic addition reaction (hidrohalogenation, hydration, halogenation, Diels-Alder, reduction and oxidation). Conjugated Addition in dienes (thermodynamic and kinetic product).6 Background of NMR, InfraRed, UV and Fluorescence techniques7.Aromatic compounds - Physical properties. Aromatic Eletrophilic Substitution . Effect of Substituent Groups. Aromatic Nucleophilic Substitution.8.Alcohols and ethers - physical properties, reactions and mechanisms."

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "210064 - Eduardo Rezende Triboni"
$ws.Range("C18").Value = "210064 - Eduardo Rezende Triboni"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Duas provas teóricas e ao longo do semestre letivoAos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."
$ws.Range("C19").Value = "Duas provas teóricas e ao longo do semestre letivoAos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A média final (M) será calculada pela expressão: M = (P1 + P2)/2"
$ws.Range("C20").Value = "A média final (M) será calculada pela expressão: M = (P1 + P2)/2"

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."
$ws.Range("C21").Value = "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."

$ws.Range("A22").Value = "Requisitos:"

$ws.Range("B23").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"

$ws.Range("B24").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"

# Adjust row heights that changed (match target customHeight / default)
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

$ws.Range("A1").Select()
